$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Mapping of row -> (DAMSLTag, DialogAct) updates, derived from the commit diff
$updates = @(
    @{ Row = 2; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 9; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 10; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 11; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 13; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 14; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 15; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 16; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 19; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 20; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 33; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 34; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 40; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 48; DAMSLTag = 'qy'; DialogAct = 'Yes-No-Question' },
    @{ Row = 53; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 55; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 68; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 72; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 73; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 82; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 84; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 90; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 99; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 100; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 106; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 111; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 112; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 113; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 123; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 145; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 149; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 170; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 175; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 181; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 183; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 185; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 186; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 203; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 206; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 207; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 208; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 210; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 228; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 234; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 238; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 240; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 244; DAMSLTag = 'b'; DialogAct = 'Acknowledge (Backchannel)' },
    @{ Row = 248; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 249; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 250; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 251; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 252; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 258; DAMSLTag = '%'; DialogAct = 'Uninterpretable' },
    @{ Row = 273; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 277; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 278; DAMSLTag = 'qy'; DialogAct = 'Yes-No-Question' },
    @{ Row = 283; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 302; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 314; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 325; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 335; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 336; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 337; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 340; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 352; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 360; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 361; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 363; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 364; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 365; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 381; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 386; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 390; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 392; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 415; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 416; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 419; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 420; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 424; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 425; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 432; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 433; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 438; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 439; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 443; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 451; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 452; DAMSLTag = 'aa'; DialogAct = 'Agree/Accept' },
    @{ Row = 457; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 460; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 469; DAMSLTag = 'ba'; DialogAct = 'Appreciation' },
    @{ Row = 497; DAMSLTag = 'sv'; DialogAct = 'Statement-opinion' },
    @{ Row = 499; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' },
    @{ Row = 502; DAMSLTag = 'sd'; DialogAct = 'Statement-non-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.DAMSLTag
    $ws.Cells.Item($u.Row, 10).Value = $u.DialogAct
}
